$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 26 ("10단위" keyword row), shifting all following rows up.
$ws.Rows.Item(26).Delete()

# Update the selection to match the post-edit state (whole row 26 selected).
$ws.Range("A26:XFD26").Select()
